$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at position 27, pushing existing rows 27+ down to 29+
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# Populate new row 27 ("Plate (part)") and row 28 ("Plate Dimensions").
# Column order chosen to reproduce the shared-string insertion order of the
# original authoring session (A27, A28, C27, D27, D28, F27, then the rest).
$ws.Range("A27").Value = "Plate (part)"
$ws.Range("A28").Value = "Plate Dimensions"
$ws.Range("C27").Value = "PhyType = 'Plate Dimensions'"
$ws.Range("D27").Value = "object/[irn]/plate-[sequence #]"
$ws.Range("D28").Value = "object/[irn]/plate-[sequence #]/[dimension type]"
$ws.Range("F27").Value = "http://vocab.getty.edu/aat/300404443"

$ws.Range("B27").Value = "X"
$ws.Range("E27").Value = "E22_Human-Made_Object"
$ws.Range("H27").Value = "https://linked.art/model/object/physical/"

$ws.Range("B28").Value = "X"
$ws.Range("E28").Value = "E54_Dimension"
$ws.Range("F28").Value = "http://vocab.getty.edu/aat/300055644`nhttp://vocab.getty.edu/aat/300055647`nhttp://vocab.getty.edu/aat/300072633`nhttp://vocab.getty.edu/aat/300055624"
$ws.Range("H28").Value = "https://linked.art/model/object/physical/"

# Row 28 displays at a taller height (matches the "Framed/Mount/..." pattern rows)
$ws.Rows.Item(28).RowHeight = 60

# Rebuild the hyperlinks collection: deleting via any hyperlink's Range wipes
# the whole sheet collection in this engine, so recreate all of them in the
# original order with addresses shifted down by 2 for rows that were >= 27.
$ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "https://linked.art/model/object/identity/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3:H9"), "https://linked.art/model/object/identity/", "", "", "https://linked.art/model/object/identity/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H13"), "https://linked.art/model/object/physical/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H36"), "https://linked.art/model/object/aboutness/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H37:H41"), "https://linked.art/model/object/aboutness/", "", "", "https://linked.art/model/object/aboutness/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H52"), "https://linked.art/model/object/provenance/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H53"), "https://linked.art/model/object/rights/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H54:H55"), "https://linked.art/model/object/rights/", "", "", "https://linked.art/model/object/rights/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H49"), "https://linked.art/model/object/provenance/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H51"), "https://linked.art/model/object/provenance/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H12"), "https://linked.art/model/base/", "statements-about-a-resource", "", "https://linked.art/model/base/ - statements-about-a-resource") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H42"), "https://linked.art/model/provenance/production.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H43"), "https://linked.art/model/provenance/production.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H45"), "https://linked.art/model/provenance/production.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H48"), "https://linked.art/model/provenance/production.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H30"), "https://linked.art/model/object/physical/", "dimensions", "", "https://linked.art/model/object/physical/ - dimensions") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I30"), "https://github.com/linked-art/linked.art/issues/191") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H56"), "https://linked.art/model/object/digital/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H57"), "https://linked.art/model/object/digital/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D57"), "http://collection.imamuseum.org/artwork/%5bdagwood%20id%5d/", "", "", "http://collection.imamuseum.org/artwork/[dagwood id]/") | Out-Null

# Restore view/selection state (top-left scroll of frozen panes isn't
# serialized by this engine, but the active cell/selection is).
$ws.Activate()
$ws.Range("H28").Select() | Out-Null
